# add - management command - generate excel dummy transactions
#
# Rebuilds the PurchaseReceiveHeader seed sheet with a full set of dummy
# rows (as if generated by a management command), using upper-cased,
# hyphenated codes (P-REC-#, P-REQ-#, TEST-ID-##) instead of the old
# "P REQ 01" style placeholders.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (overwrites the old seed header in place)
$ws.Range("A1").Value = "CODE"
$ws.Range("B1").Value = "DATE"
$ws.Range("C1").Value = "RECEIVER"
$ws.Range("D1").Value = "HEADER CODE"
$ws.Range("E1").Value = "STATUS"

# Dummy transaction rows
$rows = @(
    @{ Code = "P-REC-1";  Date = "2025-01-01"; Receiver = "TEST-ID-14"; Header = "P-REQ-1";  Status = "CLOSED" },
    @{ Code = "P-REC-2";  Date = "2025-01-02"; Receiver = "TEST-ID-15"; Header = "P-REQ-2";  Status = "CLOSED" },
    @{ Code = "P-REC-3";  Date = "2025-01-03"; Receiver = "TEST-ID-14"; Header = "P-REQ-3";  Status = "CLOSED" },
    @{ Code = "P-REC-4";  Date = "2025-01-04"; Receiver = "TEST-ID-15"; Header = "P-REQ-4";  Status = "CLOSED" },
    @{ Code = "P-REC-5";  Date = "2025-01-10"; Receiver = "TEST-ID-15"; Header = "P-REQ-5";  Status = "CLOSED" },
    @{ Code = "P-REC-6";  Date = "2025-01-14"; Receiver = "TEST-ID-14"; Header = "P-REQ-6";  Status = "CLOSED" },
    @{ Code = "P-REC-7";  Date = "2025-01-16"; Receiver = "TEST-ID-14"; Header = "P-REQ-7";  Status = "CLOSED" },
    @{ Code = "P-REC-8";  Date = "2025-01-22"; Receiver = "TEST-ID-14"; Header = "P-REQ-8";  Status = "CLOSED" },
    @{ Code = "P-REC-9";  Date = "2025-01-28"; Receiver = "TEST-ID-14"; Header = "P-REQ-9";  Status = "CLOSED" },
    @{ Code = "P-REC-10"; Date = "2025-01-31"; Receiver = "TEST-ID-14"; Header = "P-REQ-10"; Status = "CLOSED" }
)

# Apply the date format to the whole receive-date column up front, then
# fill in the values — this keeps the column on a single custom
# "yyyy-mm-dd" number format record instead of Excel minting a transient
# short-date style per write.
$ws.Range("B2:B11").NumberFormat = "yyyy\-mm\-dd"

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.Code
    $ws.Cells.Item($r, 2).Value = $row.Date
    $ws.Cells.Item($r, 3).Value = $row.Receiver
    $ws.Cells.Item($r, 4).Value = $row.Header
    $ws.Cells.Item($r, 5).Value = $row.Status
    $r++
}

$ws.Range("A2:E11").Select()
